$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: columns D (Price) and E (Volume(1h)) are stored as text in the
# source data. Values that look like plain decimals (e.g. "676.55") would
# otherwise be auto-converted to numbers by Excel, so for those we force
# the cell to Text format before writing, then restore the cell style so
# no extra formatting is left behind.

$ws.Range("D2").Value = "69.585.80"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "3.694.45"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "676.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "3.694.04"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "69.597.41"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "469.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "3.842.40"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +3.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").Value = "3.685.39"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.25%  "
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0901"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("E51").Value = "  +1.90%  "

